# Update the answer table cells to the new set of generated problems.
# The table has 20 rows (5 data rows with actual content, interleaved
# with 3 blank spacer rows each); columns are always 5.
# Data rows (1-indexed): 1, 5, 9, 13, 17.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$values = @{
    "1,1" = "77÷4=19, 1"
    "1,2" = "84÷4=21, 0"
    "1,3" = "58÷5=11, 3"
    "1,4" = "93÷4=23, 1"
    "1,5" = "58÷2=29, 0"

    "5,1" = "35÷8=4, 3"
    "5,2" = "70÷3=23, 1"
    "5,3" = "70÷7=10, 0"
    "5,4" = "77÷5=15, 2"
    "5,5" = "14÷4=3, 2"

    "9,1" = "95÷8=11, 7"
    "9,2" = "61÷5=12, 1"
    "9,3" = "70÷3=23, 1"
    "9,4" = "85÷7=12, 1"
    "9,5" = "86÷6=14, 2"

    "13,1" = "98÷2=49, 0"
    "13,2" = "60÷7=8, 4"
    "13,3" = "73÷4=18, 1"
    "13,4" = "55÷6=9, 1"
    "13,5" = "22÷4=5, 2"

    "17,1" = "39÷4=9, 3"
    "17,2" = "21÷7=3, 0"
    "17,3" = "93÷2=46, 1"
    "17,4" = "76÷3=25, 1"
    "17,5" = "44÷8=5, 4"
}

$rows = @(1, 5, 9, 13, 17)
foreach ($r in $rows) {
    for ($c = 1; $c -le 5; $c++) {
        $key = "$r,$c"
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$key]
    }
}
